# "Generate Report for Handoff" — refresh the localization-status report:
# the first file (67c54cca...) has been fully handed back and is replaced
# in the report by a newly-queued file (dc50655f...), and the second file
# (a97d0c6a...) moves from "Handed back" into "Ready for handoff" under a
# new name (ffffbbf77c52...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-09-07 07:28:18"

$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-09-07 07:28:18"

# Rebuild hyperlinks on column B (also updates A/B text + display name).
$ws1.Hyperlinks.Delete()
$ws1.Range("A2").Value = "dc50655f-ca93-4596-8aeb-64a92162e418.md"
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/67c54cca-20e7-401f-a0e0-c8583438eec2.md", "", "", "e2e\dc50655f-ca93-4596-8aeb-64a92162e418.md") | Out-Null
$ws1.Range("A3").Value = "ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md"
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/a97d0c6a-328d-4e75-a747-e0be0fa92585.md", "", "", "e2e\ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md") | Out-Null

$ws1.Columns.Item(5).ColumnWidth = 17.22
$ws1.Columns.Item(6).ColumnWidth = 17.22

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("G2").Value = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-07 07:28:11"
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Range("F3").Value = "True"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-09-07 07:28:11"
$ws2.Range("K3").Value = "0001-01-01 00:00:00"

# Target/Handback columns are now empty (no handback has happened yet).
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""

# Rebuild hyperlinks: only column A keeps a link now (I2/I3 links removed).
$ws2.Hyperlinks.Delete()
$ws2.Range("A2").Value = "dc50655f-ca93-4596-8aeb-64a92162e418.md"
$ws2.Range("A3").Value = "ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/67c54cca-20e7-401f-a0e0-c8583438eec2.md", "", "", "dc50655f-ca93-4596-8aeb-64a92162e418.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/a97d0c6a-328d-4e75-a747-e0be0fa92585.md", "", "", "ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md") | Out-Null

$ws2.Range("I2").Style = "Normal"
$ws2.Range("I3").Style = "Normal"

$ws2.Columns.Item(3).ColumnWidth = 17.22
$ws2.Columns.Item(9).ColumnWidth = 18.65
$ws2.Columns.Item(10).ColumnWidth = 21.7

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("G2").Value = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-07 07:28:18"
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Range("F3").Value = "True"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.de-de.xlf"
$ws3.Range("H3").Value = "2016-09-07 07:28:18"
$ws3.Range("K3").Value = "0001-01-01 00:00:00"

$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""

$ws3.Hyperlinks.Delete()
$ws3.Range("A2").Value = "dc50655f-ca93-4596-8aeb-64a92162e418.md"
$ws3.Range("A3").Value = "ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/67c54cca-20e7-401f-a0e0-c8583438eec2.md", "", "", "dc50655f-ca93-4596-8aeb-64a92162e418.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/a97d0c6a-328d-4e75-a747-e0be0fa92585.md", "", "", "ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md") | Out-Null

$ws3.Range("I2").Style = "Normal"
$ws3.Range("I3").Style = "Normal"

$ws3.Columns.Item(3).ColumnWidth = 17.22
$ws3.Columns.Item(9).ColumnWidth = 18.65
$ws3.Columns.Item(10).ColumnWidth = 21.7
